$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures scraped by the GitHub Actions job.
# Column D (Price) values are forced to Text format before assignment so that
# numeric-looking strings (e.g. "329.85", "1.0000") are preserved verbatim as
# text instead of being normalised/rounded as floating point numbers.

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.109.40'
$ws.Cells.Item(2, 5).Value = '  +5.57%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.918.78'
$ws.Cells.Item(3, 5).Value = '  +2.48%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.66%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '329.85'
$ws.Cells.Item(5, 5).Value = '  +4.59%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5214'
$ws.Cells.Item(7, 5).Value = '  +2.53%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4081'
$ws.Cells.Item(8, 5).Value = '  +4.68%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08564'
$ws.Cells.Item(9, 5).Value = '  +2.44%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '42.97'
$ws.Cells.Item(10, 5).Value = '  +2.66%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.78%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '22.31'
$ws.Cells.Item(12, 5).Value = '  +9.23%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.426'
$ws.Cells.Item(13, 5).Value = '  +3.42%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.927.35'
$ws.Cells.Item(14, 5).Value = '  +2.94%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.406'
$ws.Cells.Item(15, 5).Value = '  +1.86%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.64%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '95.31'
$ws.Cells.Item(17, 5).Value = '  +4.52%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +1.17%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06692'
$ws.Cells.Item(19, 5).Value = '  -0.64%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.43'

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.0000'
$ws.Cells.Item(21, 5).Value = '  -0.65%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.017'
$ws.Cells.Item(22, 5).Value = '  +1.54%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '30.123.97'
$ws.Cells.Item(23, 5).Value = '  +5.55%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.38'
$ws.Cells.Item(24, 5).Value = '  +2.54%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.201'
$ws.Cells.Item(25, 5).Value = '  +0.45%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.147.75'
$ws.Cells.Item(26, 5).Value = '  +3.08%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '21.08'
$ws.Cells.Item(27, 5).Value = '  +2.25%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '159.77'
$ws.Cells.Item(28, 5).Value = '  +0.73%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.444'
$ws.Cells.Item(29, 5).Value = '  +0.66%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '128.99'
$ws.Cells.Item(30, 5).Value = '  +2.18%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.081'
$ws.Cells.Item(31, 5).Value = '  +3.48%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +2.18%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.060'
$ws.Cells.Item(33, 5).Value = '  +5.69%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.64%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.02488'
$ws.Cells.Item(35, 5).Value = '  +1.18%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.06613'
$ws.Cells.Item(36, 5).Value = '  +0.25%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.2209'
$ws.Cells.Item(37, 5).Value = '  +2.06%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +4.25%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +2.72%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -0.18%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6547'

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.248'
$ws.Cells.Item(42, 5).Value = '  +0.59%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '11.64'
$ws.Cells.Item(43, 5).Value = '  +4.82%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6154'
$ws.Cells.Item(44, 5).Value = '  +2.49%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.20'
$ws.Cells.Item(45, 5).Value = '  +1.01%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.760'
$ws.Cells.Item(46, 5).Value = '  +2.13%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.080'
$ws.Cells.Item(47, 5).Value = '  +3.58%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.248'
$ws.Cells.Item(48, 5).Value = '  +2.80%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '124.34'
$ws.Cells.Item(49, 5).Value = '  +1.43%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.169'
$ws.Cells.Item(50, 5).Value = '  +9.79%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '79.77'

